$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(7).Insert()

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(7, 3).Value = "Maule"
$ws.Cells.Item(7, 4).Value = 44526
$ws.Cells.Item(7, 5).Value = 7
$ws.Cells.Item(7, 6).Value = "Fruta"
$ws.Cells.Item(7, 7).Value = 100101
$ws.Cells.Item(7, 8).Value = "Berries"
$ws.Cells.Item(7, 9).Value = 100101001
$ws.Cells.Item(7, 10).Value = "Arándano (blue)"
$ws.Cells.Item(7, 11).Value = "Sin especificar"
$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 250
$ws.Cells.Item(7, 14).Value = 4000
$ws.Cells.Item(7, 15).Value = 4000
$ws.Cells.Item(7, 16).Value = 4000
$ws.Cells.Item(7, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(7, 18).Value = "Provincia de Linares"
$ws.Cells.Item(7, 19).Value = 2000
$ws.Cells.Item(7, 20).Value = 2
